# Ticket 38 - add a "LessGreater" query example sheet (first_name <> 'Randy')
# and refresh sheet selections to match the authored workbook state.

$wb = $excel.ActiveWorkbook

$queryWs    = $wb.Worksheets.Item("Query")
$preparedWs = $wb.Worksheets.Item("Prepared")

# Update the selection on the existing sheets.
$queryWs.Range("A1:G2").Select() | Out-Null
$preparedWs.Range("A3").Select() | Out-Null

# Add the new "LessGreater" sheet as a copy of "Query" (same layout/styles),
# placed after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$queryWs.Copy([System.Reflection.Missing]::Value, $lastSheet) | Out-Null

$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "LessGreater"

# Swap the query text in A2 for the "not equal" variant.
$newWs.Range("A2").Value = '<jt:forEach items="${jdbc.execQuery(''SELECT * FROM employee WHERE first_name <> \''Randy\'''')}" var="employee" >${employee.first_name}'

# Make the new sheet the active tab.
$newWs.Activate() | Out-Null
$newWs.Range("A1").Select() | Out-Null
